$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1251.0714
$ws.Range("I19").Value = 885.7143
$ws.Range("J19").Value = 1372.8572
$ws.Range("K19").Value = 885.7143
$ws.Range("L19").Value = 1372.8572
$ws.Range("M19").Value = -710.7143
$ws.Range("N19").Value = -1722.8572

$ws.Range("H40").Value = 2572.7083
$ws.Range("I40").Value = 2493.4
$ws.Range("J40").Value = 2704.889
$ws.Range("K40").Value = 2493.4
$ws.Range("L40").Value = 2704.889
$ws.Range("M40").Value = -2318.4
$ws.Range("N40").Value = -3054.889

$ws.Range("H64").Value = 145214
$ws.Range("J64").Value = 2779.6
$ws.Range("L64").Value = 2779.6
$ws.Range("N64").Value = -3275.6

$ws.Range("H67").Value = 145214
$ws.Range("J67").Value = 2779.6
$ws.Range("L67").Value = 2779.6
$ws.Range("N67").Value = -4495.6

$ws.Range("H106").Value = 2219.9167
$ws.Range("I106").Value = 2079.875
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 2079.875
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -1448.875
$ws.Range("N106").Value = -3762

$ws.Range("H107").Value = 476.47827
$ws.Range("I107").Value = 438.41177
$ws.Range("J107").Value = 584.3333
$ws.Range("K107").Value = 438.41177
$ws.Range("L107").Value = 584.3333
$ws.Range("M107").Value = 1481.58823
$ws.Range("N107").Value = -4424.3333

$ws.Range("H112").Value = 983
$ws.Range("I112").Value = 650
$ws.Range("J112").Value = 1009.64
$ws.Range("K112").Value = 1950
$ws.Range("L112").Value = 3028.92
$ws.Range("M112").Value = -842
$ws.Range("N112").Value = -5244.92

$ws.Range("H129").Value = 3121.8667
$ws.Range("I129").Value = 8195
$ws.Range("J129").Value = 1060.9062
$ws.Range("K129").Value = 24585
$ws.Range("L129").Value = 3182.7186
$ws.Range("M129").Value = -19585
$ws.Range("N129").Value = -13182.7186

$ws.Range("H132").Value = 4390494
$ws.Range("I132").Value = 4721450
$ws.Range("J132").Value = 5327.75
$ws.Range("K132").Value = 14164350
$ws.Range("L132").Value = 15983.25
$ws.Range("M132").Value = -14161820
$ws.Range("N132").Value = -21043.25

$ws.Range("H135").Value = 503.78333
$ws.Range("I135").Value = 453.32074
$ws.Range("J135").Value = 885.8570999999999
$ws.Range("K135").Value = 4079.88666
$ws.Range("L135").Value = 7972.7139
$ws.Range("M135").Value = -1544.88666
$ws.Range("N135").Value = -13042.7139

$ws.Range("H137").Value = 1904.1818
$ws.Range("I137").Value = 1581.9412
$ws.Range("J137").Value = 2999.8
$ws.Range("K137").Value = 4745.8236
$ws.Range("L137").Value = 8999.400000000001
$ws.Range("M137").Value = -2195.8236
$ws.Range("N137").Value = -14099.4

$ws.Range("H138").Value = 1826.2354
$ws.Range("I138").Value = 1546.5927
$ws.Range("J138").Value = 2904.8572
$ws.Range("K138").Value = 4639.7781
$ws.Range("L138").Value = 8714.571599999999
$ws.Range("M138").Value = 500.2219000000005
$ws.Range("N138").Value = -18994.5716

$ws.Range("H141").Value = 2001.6444
$ws.Range("I141").Value = 1842.6923
$ws.Range("J141").Value = 3034.8333
$ws.Range("K141").Value = 5528.0769
$ws.Range("L141").Value = 9104.499899999999
$ws.Range("M141").Value = -348.0769
$ws.Range("N141").Value = -19464.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 5688.5
$ws.Range("J43").Value = 5688.5
$ws.Range("L43").Value = 5688.5
$ws.Range("N43").Value = -6314.5

$ws.Range("H61").Value = 2461.7942
$ws.Range("I61").Value = 1409.7
$ws.Range("J61").Value = 2900.1667
$ws.Range("K61").Value = 1409.7
$ws.Range("L61").Value = 2900.1667
$ws.Range("M61").Value = -1197.7
$ws.Range("N61").Value = -3324.1667

$ws.Range("H74").Value = 708.5789
$ws.Range("I74").Value = 704.2
$ws.Range("J74").Value = 725
$ws.Range("K74").Value = 704.2
$ws.Range("L74").Value = 725
$ws.Range("M74").Value = 169.8
$ws.Range("N74").Value = -2473

$ws.Range("H77").Value = 708.5789
$ws.Range("I77").Value = 704.2
$ws.Range("J77").Value = 725
$ws.Range("K77").Value = 3521
$ws.Range("L77").Value = 3625
$ws.Range("M77").Value = 847
$ws.Range("N77").Value = -12361

$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

$ws.Range("H136").Value = 2461.7942
$ws.Range("I136").Value = 1409.7
$ws.Range("J136").Value = 2900.1667
$ws.Range("K136").Value = 4229.1
$ws.Range("L136").Value = 8700.500100000001
$ws.Range("M136").Value = -1679.1
$ws.Range("N136").Value = -13800.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1660.9744
$ws.Range("I99").Value = 1194.5
$ws.Range("J99").Value = 2060.8096
$ws.Range("K99").Value = 1194.5
$ws.Range("L99").Value = 2060.8096
$ws.Range("M99").Value = 303.5
$ws.Range("N99").Value = -5056.809600000001

$ws.Range("H107").Value = 45473176
$ws.Range("I107").Value = 55578156
$ws.Range("J107").Value = 776.25
$ws.Range("K107").Value = 55578156
$ws.Range("L107").Value = 776.25
$ws.Range("M107").Value = -55576236
$ws.Range("N107").Value = -4616.25

$ws.Range("H134").Value = 4459.108
$ws.Range("I134").Value = 4499
$ws.Range("J134").Value = 4364.8184
$ws.Range("K134").Value = 13497
$ws.Range("L134").Value = 13094.4552
$ws.Range("M134").Value = -10962
$ws.Range("N134").Value = -18164.4552

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000600
$ws.Range("I4").Value = 5000600
$ws.Range("K4").Value = 5000600
$ws.Range("M4").Value = -5000488

$ws.Range("H22").Value = 442.25
$ws.Range("I22").Value = 442.25
$ws.Range("K22").Value = 442.25
$ws.Range("M22").Value = -92.25

$ws.Range("H31").Value = 27104.518
$ws.Range("I31").Value = 1247.4231
$ws.Range("K31").Value = 1247.4231
$ws.Range("M31").Value = -952.4231

$ws.Range("H34").Value = 27104.518
$ws.Range("I34").Value = 1247.4231
$ws.Range("K34").Value = 1247.4231
$ws.Range("M34").Value = -1045.4231

$ws.Range("H86").Value = 2985.6924
$ws.Range("J86").Value = 2951.5
$ws.Range("L86").Value = 2951.5
$ws.Range("N86").Value = -5197.5

$ws.Range("H89").Value = 2985.6924
$ws.Range("J89").Value = 2951.5
$ws.Range("L89").Value = 14757.5
$ws.Range("N89").Value = -25989.5

$ws.Range("H107").Value = 919.8823
$ws.Range("I107").Value = 1112.5
$ws.Range("K107").Value = 1112.5
$ws.Range("M107").Value = 807.5

$ws.Range("H132").Value = 27275504
$ws.Range("I132").Value = 23812178
$ws.Range("J132").Value = 38464708
$ws.Range("K132").Value = 71436534
$ws.Range("L132").Value = 115394124
$ws.Range("M132").Value = -71434004
$ws.Range("N132").Value = -115399184

$ws.Range("H134").Value = 1023.9796
$ws.Range("I134").Value = 886.8182
$ws.Range("K134").Value = 2660.4546
$ws.Range("M134").Value = -125.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2531.111
$ws.Range("J109").Value = 3120
$ws.Range("L109").Value = 9360
$ws.Range("N109").Value = -11440

$ws.Range("H129").Value = 17866316
$ws.Range("I129").Value = 83333710
$ws.Range("J129").Value = 11570.909
$ws.Range("K129").Value = 250001130
$ws.Range("L129").Value = 34712.727
$ws.Range("M129").Value = -249996130
$ws.Range("N129").Value = -44712.727

$ws.Range("H131").Value = 1484
$ws.Range("J131").Value = 1502.5177
$ws.Range("L131").Value = 4507.5531
$ws.Range("N131").Value = -14587.5531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 130120.125
$ws.Range("I70").Value = 253739.5
$ws.Range("J70").Value = 6500.75
$ws.Range("K70").Value = 253739.5
$ws.Range("L70").Value = 6500.75
$ws.Range("M70").Value = -253469.5
$ws.Range("N70").Value = -7040.75

$ws.Range("H73").Value = 130120.125
$ws.Range("I73").Value = 253739.5
$ws.Range("J73").Value = 6500.75
$ws.Range("K73").Value = 253739.5
$ws.Range("L73").Value = 6500.75
$ws.Range("M73").Value = -252803.5
$ws.Range("N73").Value = -8372.75

$ws.Range("H107").Value = 842124.5
$ws.Range("I107").Value = 407.63635
$ws.Range("K107").Value = 407.63635
$ws.Range("M107").Value = 1512.36365

$ws.Range("H111").Value = 17125
$ws.Range("J111").Value = 17125
$ws.Range("L111").Value = 17125
$ws.Range("N111").Value = -23259

$ws.Range("H113").Value = 1720.25
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 1823.2858
$ws.Range("K113").Value = 999
$ws.Range("L113").Value = 1823.2858
$ws.Range("M113").Value = 1171
$ws.Range("N113").Value = -6163.2858

$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws.Range("H132").Value = 3694
$ws.Range("I132").Value = 3507.25
$ws.Range("J132").Value = 4067.5
$ws.Range("K132").Value = 10521.75
$ws.Range("L132").Value = 12202.5
$ws.Range("M132").Value = -7991.75
$ws.Range("N132").Value = -17262.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws.Range("H120").Value = 40333.668
$ws.Range("J120").Value = 40333.668
$ws.Range("L120").Value = 40333.668
$ws.Range("N120").Value = -50009.668

$ws.Range("H132").Value = 3529.8572
$ws.Range("I132").Value = 3383.56
$ws.Range("K132").Value = 10150.68
$ws.Range("M132").Value = -7620.68

$ws.Range("H136").Value = 1860.2667
$ws.Range("I136").Value = 1719.9
$ws.Range("J136").Value = 2141
$ws.Range("K136").Value = 5159.700000000001
$ws.Range("L136").Value = 6423
$ws.Range("M136").Value = -2609.700000000001
$ws.Range("N136").Value = -11523

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2524000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 2524000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 2524000
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = -2524224

$ws.Range("H132").Value = 2270.6606
$ws.Range("I132").Value = 2384.262
$ws.Range("J132").Value = 1929.8572
$ws.Range("K132").Value = 7152.786
$ws.Range("L132").Value = 5789.571599999999
$ws.Range("M132").Value = -4622.786
$ws.Range("N132").Value = -10849.5716

$ws.Range("H135").Value = 40312.582
$ws.Range("J135").Value = 40312.582
$ws.Range("L135").Value = 40312.582
$ws.Range("N135").Value = -50452.582
